$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (second table, Tableau13: I=ID, J=Taches, K=Commentaires, L=Visa, M=Responsable1, N=Responsable2, O=Etat)
# J7: new task title (keeps existing style s=3)
$ws.Range("J7").Value = "HOF - Total bug game"

# K7: bug description / comment (keeps existing style s=4)
$ws.Range("K7").Value = "Quand tu vas dans le hall of fame et que tu select un niveau et tu appuies sur la flèche du bas plantage total"

# L7: visa date -> 2020-06-01 (serial 43983), style changes to the date-formatted style (like L5 / D3 / D4)
$ws.Range("L5").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L7").Value = 43983

# M7: responsable -> "Brice" (keeps existing style s=4)
$ws.Range("M7").Value = "Brice"

# O7: state -> "Terminé (...)" with the "Terminé" green italic style (same as O5/O6)
$ws.Range("O6").Copy()
$ws.Range("O7").PasteSpecial(-4122)
$ws.Range("O7").Value = "Terminé (prevent default ligne 54 menu.js)"

# Column O (15) widens slightly to fit the new longest "Terminé" text
$ws.Columns("O:O").ColumnWidth = 38.2

# Update the active selection to J4 (matches the recorded UI state after editing)
$ws.Range("J4").Select() | Out-Null

$excel.CutCopyMode = $false
